$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) "Mes motivations" intro sentence: insert "une histoire," so that
#    "... vous racontez mon histoire et ce qui ..."
#    becomes
#    "... vous racontez une histoire, mon histoire et ce qui ..."
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "vous racontez mon histoire et ce qui", $true, $false, $false, $false, $false,
    $true, 1, $false, "vous racontez une histoire, mon histoire et ce qui", 2) | Out-Null

# -----------------------------------------------------------------
# 2) Drop the trailing stray space that used to sit at the end of the
#    "Pour réaliser cette technique ..." paragraph (the one ending in
#    "... en question). ") - it is trimmed to "... en question)."
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "en question). ", $true, $false, $false, $false, $false,
    $true, 1, $false, "en question).", 2) | Out-Null

# -----------------------------------------------------------------
# 3) Relocate the (otherwise invisible) "_GoBack" bookmark: it used to
#    sit between "... sur le" and " corps ..." - it now belongs a
#    couple of paragraphs later, right after "... (qui sont relié" and
#    before " au bone du haut du corps) ..." near "chaque bras".
# -----------------------------------------------------------------
$target = $d.Content.Duplicate
$target.Find.Execute(
    "puis deux pour chaque bras (qui sont relié", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$bookmarkPoint = $d.Range($target.End, $target.End)
$d.Bookmarks.Add("_GoBack", $bookmarkPoint) | Out-Null
